$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Phase 1: write every new cell VALUE in the same order the strings were
# first introduced (keeps the shared-string table ordering faithful).
# ---------------------------------------------------------------------
$ws.Range("C7").Value  = "5eaf693"
$ws.Range("B9").Value  = "Product Definition"
$ws.Range("D9").Value  = "Not Uploaded"
$ws.Range("E9").Value  = "Complete it and upload it"
$ws.Range("B10").Value = "Flow Daigram"
$ws.Range("D10").Value = "See First Draft_Block diagram - MK Comments for details"
$ws.Range("E10").Value = "This document does not match the guidelines"
$ws.Range("B3").Value  = "Arduino Electronic Load "
$ws.Range("B7").Value  = "File Struture"
$ws.Range("D7").Value  = "Files from previous design clutter up the folder structure"
$ws.Range("E7").Value  = "As you add you own files to the folders move the files from the previous design to the Archive folder. This makes it clear which are the active files and make them easier to find"
$ws.Range("D8").Value  = "Version in file names are not needed because Github hands that"
$ws.Range("E8").Value  = "Remove the version from the file name"

$ws.Range("A7").Value  = "Malcolm Knapp"
$ws.Range("A8").Value  = "Malcolm Knapp"
$ws.Range("A9").Value  = "Malcolm Knapp"
$ws.Range("A10").Value = "Malcolm Knapp"
$ws.Range("B8").Value  = "File Struture"
$ws.Range("C8").Value  = "5eaf693"
$ws.Range("C9").Value  = "5eaf693"
$ws.Range("C10").Value = "5eaf693"

# ---------------------------------------------------------------------
# Phase 2: formatting
# ---------------------------------------------------------------------

# Column A (Reviewer) on the new rows - bold Arial, wrapped (matches the
# header look already used for A6/F6).
foreach ($addr in "A7","A8","A9","A10") {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Name = "Arial"
    $c.Font.Size = 13
    $c.Font.Color = 0
    $c.WrapText = $true
}

# Column B (Document) on rows 7 & 8 - Calibri 13, black, wrapped.
foreach ($addr in "B7","B8") {
    $c = $ws.Range($addr)
    $c.WrapText = $true
    $c.Font.Name = "Calibri"
    $c.Font.Size = 13
    $c.Font.Color = 0
}

# Column C (Version / commit hash) on all 4 new rows - Segoe UI 12, dark gray.
foreach ($addr in "C7","C8","C9","C10") {
    $c = $ws.Range($addr)
    $c.Font.Name = "Segoe UI"
    $c.Font.Size = 12
    $c.Font.Color = 3025188
}

# Columns D & E (Issue / Suggestion) on rows 7 & 8 - Arial 13, black, wrapped.
foreach ($addr in "D7","E7","D8","E8") {
    $c = $ws.Range($addr)
    $c.WrapText = $true
    $c.Font.Name = "Arial"
    $c.Font.Size = 13
    $c.Font.Color = 0
}

# Trailing empty cells on rows 7/8 styled like the rest of column A.
foreach ($addr in "F7","F8") {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Name = "Arial"
    $c.Font.Size = 13
    $c.Font.Color = 0
    $c.WrapText = $true
}

# Row heights for the two wrapped feedback rows, and the Project row.
$ws.Rows.Item(3).RowHeight = 32
$ws.Rows.Item(7).RowHeight = 48
$ws.Rows.Item(8).RowHeight = 32

# Leave the selection where the author left it when they saved.
$ws.Range("B9").Select()

Write-Output "Feedback rows added"
